# Edit script: append processed data for 2020-11-23 .. 2020-11-26
# (commit message: "processed data for 2020-11-26")
# Extends the ALL_AGE_FINAL sheet / defined name from A1:H2479 to A1:H2519
# (4 new report dates x 10 age-range rows each = 40 new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy formatting from the last existing row (2479) down onto the
#        40 new rows so the new DATE cells (col A) pick up the existing
#        date-number-format style, matching the rest of the column. ---
$ws.Range("A2479:H2479").Copy() | Out-Null
$ws.Range("A2480:H2519").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Write the 40 new data rows (4 dates x 10 age buckets). ---
$arr = New-Object 'object[,]' 40,8
$arr[0,0] = 44158
$arr[0,1] = "0-10 years"
$arr[0,2] = 16983
$arr[0,3] = 0.0492903787548977
$arr[0,4] = 231
$arr[0,5] = 0.0567010309278351
$arr[0,6] = 4
$arr[0,7] = 0
$arr[1,0] = 44158
$arr[1,1] = "11-20 years"
$arr[1,2] = 45947
$arr[1,3] = 0.133353649687999
$arr[1,4] = 574
$arr[1,5] = 0.140893470790378
$arr[1,6] = 2
$arr[1,7] = 0
$arr[2,0] = 44158
$arr[2,1] = "21-30 years"
$arr[2,2] = 66533
$arr[2,3] = 0.19310114642287
$arr[2,4] = 709
$arr[2,5] = 0.174030436917035
$arr[2,6] = 27
$arr[2,7] = 1
$arr[3,0] = 44158
$arr[3,1] = "31-40 years"
$arr[3,2] = 54441
$arr[3,3] = 0.1580060949064
$arr[3,4] = 584
$arr[3,5] = 0.143348060873834
$arr[3,6] = 57
$arr[3,7] = 1
$arr[4,0] = 44158
$arr[4,1] = "41-50 years"
$arr[4,2] = 51173
$arr[4,3] = 0.148521259613989
$arr[4,4] = 591
$arr[4,5] = 0.145066273932253
$arr[4,6] = 147
$arr[4,7] = 1
$arr[5,0] = 44158
$arr[5,1] = "51-60 years"
$arr[5,2] = 46422
$arr[5,3] = 0.134732259468872
$arr[5,4] = 644
$arr[5,5] = 0.15807560137457
$arr[5,6] = 399
$arr[5,7] = 4
$arr[6,0] = 44158
$arr[6,1] = "61-70 years"
$arr[6,2] = 32571
$arr[6,3] = 0.0945319982585982
$arr[6,4] = 396
$arr[6,5] = 0.0972017673048601
$arr[6,6] = 818
$arr[6,7] = 5
$arr[7,0] = 44158
$arr[7,1] = "71-80 years"
$arr[7,2] = 19489
$arr[7,3] = 0.0565636337251487
$arr[7,4] = 228
$arr[7,5] = 0.0559646539027982
$arr[7,6] = 1298
$arr[7,7] = 12
$arr[8,0] = 44158
$arr[8,1] = "81+ years"
$arr[8,2] = 10482
$arr[8,3] = 0.0304222899434044
$arr[8,4] = 102
$arr[8,5] = 0.0250368188512518
$arr[8,6] = 1549
$arr[8,7] = 11
$arr[9,0] = 44158
$arr[9,1] = "Pending"
$arr[9,2] = 509
$arr[9,3] = 0.00147728921782035
$arr[9,4] = 15
$arr[9,5] = 0.00368188512518409
$arr[9,6] = 0
$arr[9,7] = 0
$arr[10,0] = 44159
$arr[10,1] = "0-10 years"
$arr[10,2] = 17049
$arr[10,3] = 0.0492953674093693
$arr[10,4] = 66
$arr[10,5] = 0.0506134969325153
$arr[10,6] = 4
$arr[10,7] = 0
$arr[11,0] = 44159
$arr[11,1] = "11-20 years"
$arr[11,2] = 46115
$arr[11,3] = 0.133336610245942
$arr[11,4] = 168
$arr[11,5] = 0.128834355828221
$arr[11,6] = 2
$arr[11,7] = 0
$arr[12,0] = 44159
$arr[12,1] = "21-30 years"
$arr[12,2] = 66742
$arr[12,3] = 0.192977383520214
$arr[12,4] = 209
$arr[12,5] = 0.160276073619632
$arr[12,6] = 27
$arr[12,7] = 0
$arr[13,0] = 44159
$arr[13,1] = "31-40 years"
$arr[13,2] = 54617
$arr[13,3] = 0.157919237597368
$arr[13,4] = 176
$arr[13,5] = 0.134969325153374
$arr[13,6] = 58
$arr[13,7] = 1
$arr[14,0] = 44159
$arr[14,1] = "41-50 years"
$arr[14,2] = 51371
$arr[14,3] = 0.148533774367219
$arr[14,4] = 198
$arr[14,5] = 0.151840490797546
$arr[14,6] = 148
$arr[14,7] = 1
$arr[15,0] = 44159
$arr[15,1] = "51-60 years"
$arr[15,2] = 46606
$arr[15,3] = 0.134756284443725
$arr[15,4] = 184
$arr[15,5] = 0.141104294478528
$arr[15,6] = 403
$arr[15,7] = 4
$arr[16,0] = 44159
$arr[16,1] = "61-70 years"
$arr[16,2] = 32723
$arr[16,3] = 0.0946150687862508
$arr[16,4] = 152
$arr[16,5] = 0.116564417177914
$arr[16,6] = 828
$arr[16,7] = 10
$arr[17,0] = 44159
$arr[17,1] = "71-80 years"
$arr[17,2] = 19601
$arr[17,3] = 0.0566742035656664
$arr[17,4] = 112
$arr[17,5] = 0.0858895705521472
$arr[17,6] = 1317
$arr[17,7] = 19
$arr[18,0] = 44159
$arr[18,1] = "81+ years"
$arr[18,2] = 10527
$arr[18,3] = 0.0304376991447258
$arr[18,4] = 45
$arr[18,5] = 0.0345092024539877
$arr[18,6] = 1587
$arr[18,7] = 38
$arr[19,0] = 44159
$arr[19,1] = "Pending"
$arr[19,2] = 503
$arr[19,3] = 0.00145437091952095
$arr[19,4] = -6
$arr[19,5] = -0.00460122699386503
$arr[19,6] = 0
$arr[19,7] = 0
$arr[20,0] = 44160
$arr[20,1] = "0-10 years"
$arr[20,2] = 17157
$arr[20,3] = 0.0493056912625154
$arr[20,4] = 108
$arr[20,5] = 0.0509915014164306
$arr[20,6] = 4
$arr[20,7] = 0
$arr[21,0] = 44160
$arr[21,1] = "11-20 years"
$arr[21,2] = 46433
$arr[21,3] = 0.133438897382548
$arr[21,4] = 318
$arr[21,5] = 0.15014164305949
$arr[21,6] = 2
$arr[21,7] = 0
$arr[22,0] = 44160
$arr[22,1] = "21-30 years"
$arr[22,2] = 67048
$arr[22,3] = 0.192682169829756
$arr[22,4] = 306
$arr[22,5] = 0.144475920679887
$arr[22,6] = 27
$arr[22,7] = 0
$arr[23,0] = 44160
$arr[23,1] = "31-40 years"
$arr[23,2] = 54878
$arr[23,3] = 0.157708091455634
$arr[23,4] = 261
$arr[23,5] = 0.123229461756374
$arr[23,6] = 58
$arr[23,7] = 0
$arr[24,0] = 44160
$arr[24,1] = "41-50 years"
$arr[24,2] = 51682
$arr[24,3] = 0.148523444415068
$arr[24,4] = 311
$arr[24,5] = 0.146836638338055
$arr[24,6] = 150
$arr[24,7] = 2
$arr[25,0] = 44160
$arr[25,1] = "51-60 years"
$arr[25,2] = 46902
$arr[25,3] = 0.134786706976423
$arr[25,4] = 296
$arr[25,5] = 0.139754485363551
$arr[25,6] = 416
$arr[25,7] = 13
$arr[26,0] = 44160
$arr[26,1] = "61-70 years"
$arr[26,2] = 32937
$arr[26,3] = 0.0946541675767016
$arr[26,4] = 214
$arr[26,5] = 0.101038715769594
$arr[26,6] = 843
$arr[26,7] = 15
$arr[27,0] = 44160
$arr[27,1] = "71-80 years"
$arr[27,2] = 19787
$arr[27,3] = 0.0568637706482131
$arr[27,4] = 186
$arr[27,5] = 0.0878186968838527
$arr[27,6] = 1339
$arr[27,7] = 22
$arr[28,0] = 44160
$arr[28,1] = "81+ years"
$arr[28,2] = 10647
$arr[28,3] = 0.0305972894370811
$arr[28,4] = 120
$arr[28,5] = 0.056657223796034
$arr[28,6] = 1626
$arr[28,7] = 39
$arr[29,0] = 44160
$arr[29,1] = "Pending"
$arr[29,2] = 501
$arr[29,3] = 0.00143977101605876
$arr[29,4] = -2
$arr[29,5] = -0.000944287063267233
$arr[29,6] = 1
$arr[29,7] = 1
$arr[30,0] = 44161
$arr[30,1] = "0-10 years"
$arr[30,2] = 17372
$arr[30,3] = 0.0492996117783277
$arr[30,4] = 215
$arr[30,5] = 0.048819255222525
$arr[30,6] = 4
$arr[30,7] = 0
$arr[31,0] = 44161
$arr[31,1] = "11-20 years"
$arr[31,2] = 47002
$arr[31,3] = 0.133385928667106
$arr[31,4] = 569
$arr[31,5] = 0.129200726612171
$arr[31,6] = 2
$arr[31,7] = 0
$arr[32,0] = 44161
$arr[32,1] = "21-30 years"
$arr[32,2] = 67824
$arr[32,3] = 0.192476218584694
$arr[32,4] = 776
$arr[32,5] = 0.176203451407811
$arr[32,6] = 28
$arr[32,7] = 1
$arr[33,0] = 44161
$arr[33,1] = "31-40 years"
$arr[33,2] = 55525
$arr[33,3] = 0.157573160487661
$arr[33,4] = 647
$arr[33,5] = 0.146911898274296
$arr[33,6] = 58
$arr[33,7] = 0
$arr[34,0] = 44161
$arr[34,1] = "41-50 years"
$arr[34,2] = 52304
$arr[34,3] = 0.14843235634663
$arr[34,4] = 622
$arr[34,5] = 0.141235240690282
$arr[34,6] = 153
$arr[34,7] = 3
$arr[35,0] = 44161
$arr[35,1] = "51-60 years"
$arr[35,2] = 47513
$arr[35,3] = 0.134836084182805
$arr[35,4] = 611
$arr[35,5] = 0.138737511353315
$arr[35,6] = 419
$arr[35,7] = 3
$arr[36,0] = 44161
$arr[36,1] = "61-70 years"
$arr[36,2] = 33428
$arr[36,3] = 0.0948645764751288
$arr[36,4] = 491
$arr[36,5] = 0.111489554950045
$arr[36,6] = 854
$arr[36,7] = 11
$arr[37,0] = 44161
$arr[37,1] = "71-80 years"
$arr[37,2] = 20076
$arr[37,3] = 0.0569732331373306
$arr[37,4] = 289
$arr[37,5] = 0.065622161671208
$arr[37,6] = 1354
$arr[37,7] = 15
$arr[38,0] = 44161
$arr[38,1] = "81+ years"
$arr[38,2] = 10817
$arr[38,3] = 0.0306973233137331
$arr[38,4] = 170
$arr[38,5] = 0.0386012715712988
$arr[38,6] = 1646
$arr[38,7] = 20
$arr[39,0] = 44161
$arr[39,1] = "Pending"
$arr[39,2] = 515
$arr[39,3] = 0.00146150702658524
$arr[39,4] = 14
$arr[39,5] = 0.00317892824704814
$arr[39,6] = 1
$arr[39,7] = 0

$ws.Range("A2480:H2519").Value = $arr

# --- 3. Update the defined name so it covers the new range. ---
$wb.Names.Item("ALL_AGE_FINAL").RefersTo = "='ALL_AGE_FINAL'!" + '$A$1:$H$2519'
